$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.775.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.60%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.446.45"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.38"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.25"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.93%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.437.38"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.31%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.19%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.10"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.889.42"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.72%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.764.98"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.452.92"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.62"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.73%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.48"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.06"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.92"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.06"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.14"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "582.96"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -7.10%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0923"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.60%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.56%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.77%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.53%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "150.85"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.37"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.14"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.20%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.72"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.37"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.81%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0291"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +25.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "143.42"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.65%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.69%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0511"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.96%  "
